# Updated cryptos list on Mon May 27 16:48:33 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking figures (e.g. "612.23", "70.346.22") that
# are stored as plain text in this sheet, not real numbers (note the
# thousands-grouped values like "70.346.22" and price figures that keep
# significant trailing zeros like "171.10"). Force text formatting before
# writing so Excel does not silently coerce these into numbers, then restore
# the default "Normal" style so the cell formatting matches the rest of the
# sheet once the text is safely stored.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.346.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.955.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.31%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.953.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("E9").Value = "  +1.51%  "

$ws.Range("E10").Value = "  +2.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.60%  "

$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.40%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.621.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.972.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.222.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("E18").Value = "  +1.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "498.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.55%  "

$ws.Range("E23").Value = "  +3.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000167"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "

$ws.Range("E26").Value = "  +3.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.14%  "

$ws.Range("E28").Value = "  +2.94%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("E30").Value = "  +1.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.110.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.920.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.70%  "

$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("E40").Value = "  +9.79%  "

$ws.Range("E41").Value = "  +3.76%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("E43").Value = "  +7.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "442.53"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("E46").Value = "  +3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +23.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0369"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.68%  "

# Row 50 is now Arweave (was Monero); row 51 is now Monero (was Arweave) -
# the two coins swapped ranking places and each carries its own refreshed
# price/volume figures.
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.81%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
